$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D holds numeric-looking text (prices such as "1.630.36" or "0.570").
# Excel auto-converts plain numeric-looking strings assigned via .Value into real
# Number cells (dropping significant trailing zeros, e.g. "0.570" -> 0.57), so for
# every column D write we snapshot/restore the cell Style around a temporary
# Text (@) NumberFormat to force the literal string to stick without altering the
# cells persisted style index.
function Set-TextValue($cell, $value) {
    $origStyle = $cell.Style
    $cell.NumberFormat = "@"
    $cell.Value = $value
    $cell.Style = $origStyle
}

Set-TextValue $ws.Range('D2') '29.873.02'
$ws.Range('E2').Value = '  +0.79%  '
Set-TextValue $ws.Range('D3') '1.630.36'
$ws.Range('E3').Value = '  +1.13%  '
$ws.Range('E4').Value = '  +0.71%  '
Set-TextValue $ws.Range('D5') '214.86'
$ws.Range('E5').Value = '  +1.15%  '
$ws.Range('E6').Value = '  +0.21%  '
Set-TextValue $ws.Range('D7') '0.999'
$ws.Range('E7').Value = '  +0.71%  '
Set-TextValue $ws.Range('D8') '28.72'
$ws.Range('E8').Value = '  -1.07%  '
$ws.Range('E9').Value = '  +0.32%  '
Set-TextValue $ws.Range('D10') '0.0609'
$ws.Range('E10').Value = '  +0.06%  '
Set-TextValue $ws.Range('D11') '0.0898'
$ws.Range('E11').Value = '  -1.21%  '
Set-TextValue $ws.Range('D12') '1.864.29'
$ws.Range('E12').Value = '  +1.12%  '
Set-TextValue $ws.Range('D13') '1.633.25'
$ws.Range('E13').Value = '  +1.25%  '
Set-TextValue $ws.Range('D14') '0.570'
$ws.Range('E14').Value = '  +0.31%  '
$ws.Range('E15').Value = '  +5.79%  '
Set-TextValue $ws.Range('D16') '29.872.78'
$ws.Range('E16').Value = '  +0.77%  '
$ws.Range('E17').Value = '  -0.97%  '
Set-TextValue $ws.Range('D18') '65.16'
$ws.Range('E18').Value = '  +1.65%  '
Set-TextValue $ws.Range('D19') '240.50'
$ws.Range('E19').Value = '  -0.16%  '
$ws.Range('E20').Value = '  -1.11%  '
$ws.Range('E21').Value = '  +0.56%  '
$ws.Range('B22').Value = 'Avalanche'
$ws.Range('C22').Value = 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
Set-TextValue $ws.Range('D22') '9.82'
$ws.Range('E22').Value = '  +1.75%  '
$ws.Range('B23').Value = 'Uniswap'
$ws.Range('C23').Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
Set-TextValue $ws.Range('D23') '4.13'
$ws.Range('E23').Value = '  +0.94%  '
$ws.Range('E24').Value = '  +3.29%  '
Set-TextValue $ws.Range('D25') '157.50'
$ws.Range('E26').Value = '  -0.80%  '
$ws.Range('E27').Value = '  -0.71%  '
$ws.Range('E28').Value = '  +0.29%  '
$ws.Range('E29').Value = '  +0.61%  '
Set-TextValue $ws.Range('D30') '0.0489'
$ws.Range('E30').Value = '  +0.49%  '
$ws.Range('E31').Value = '  +2.00%  '
$ws.Range('E32').Value = '  +2.18%  '
$ws.Range('E33').Value = '  -0.72%  '
Set-TextValue $ws.Range('D34') '1.423.97'
$ws.Range('E34').Value = '  -0.02%  '
$ws.Range('E35').Value = '  +3.35%  '
$ws.Range('E36').Value = '  -2.46%  '
Set-TextValue $ws.Range('D37') '2.76'
$ws.Range('E37').Value = '  -3.92%  '
$ws.Range('E38').Value = '  +0.57%  '
$ws.Range('E39').Value = '  +0.48%  '
$ws.Range('B40').Value = 'Aave'
$ws.Range('C40').Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
Set-TextValue $ws.Range('D40') '75.02'
$ws.Range('E40').Value = '  +7.57%  '
$ws.Range('B41').Value = 'ImmutableX'
$ws.Range('C41').Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
Set-TextValue $ws.Range('D41') '0.555'
$ws.Range('E41').Value = '  +0.06%  '
Set-TextValue $ws.Range('D42') '0.0504'
$ws.Range('E42').Value = '  -0.01%  '
$ws.Range('B43').Value = 'RenderToken'
$ws.Range('C43').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
Set-TextValue $ws.Range('D43') '1.99'
$ws.Range('E43').Value = '  +0.75%  '
$ws.Range('B44').Value = 'ARBITRUM'
$ws.Range('C44').Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
Set-TextValue $ws.Range('D44') '0.832'
$ws.Range('E44').Value = '  +0.45%  '
$ws.Range('E46').Value = '  +0.73%  '
Set-TextValue $ws.Range('D47') '1.771.50'
$ws.Range('E47').Value = '  +1.04%  '
Set-TextValue $ws.Range('D48') '5.32'
$ws.Range('E48').Value = '  -1.99%  '
Set-TextValue $ws.Range('D49') '48.56'
$ws.Range('E49').Value = '  -9.15%  '
Set-TextValue $ws.Range('D50') '92.12'
$ws.Range('E50').Value = '  +5.05%  '
$ws.Range('E51').Value = '  +1.74%  '
